$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 by shifting row 3 down to row 4.
$ws.Rows.Item(4).Insert()

# Copy the old row 3 contents (now still at row 3) into the new row 4.
$ws.Range("A3:R3").Copy()
$ws.Range("A4").PasteSpecial()

# Update row 3 with the new weekly values.
$ws.Range("D3").Value = 44848
$ws.Range("J3").Value = 1000
